$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for the
# c0354d71-... row (row 3) to reflect the newly generated handback report.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G3").Value = "2016-08-24 18:57:29"

# zh-cn sheet: update Correspond Handoff Datetime (H3) and
# Correspond Handback DateTime (K3) for the c0354d71-... row (row 3).
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("H3").Value = "2016-08-24 18:57:24"
$ws2.Range("K3").Value = "2016-08-24 18:57:40"

# de-de sheet: update Correspond Handoff Datetime (H3) and
# Correspond Handback DateTime (K3) for the c0354d71-... row (row 3).
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("H3").Value = "2016-08-24 18:57:29"
$ws3.Range("K3").Value = "2016-08-24 18:57:48"
